# Add a new "2020" column (K) to the CITES permits table, mirroring the
# formatting of the existing last data column (J), then move the active
# selection to I18 (as captured in the saved worksheet view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column J (rows 3-5) carries the per-column formatting we need to extend
# into the new column K:
#   J3 - bottom border only (header rule row)
#   J4 - bold/right-aligned "year" style
#   J5 - bordered "value" style
# Copying it over before writing values reproduces the same style indices
# (s="4" / s="5" / s="7") that a user would get by dragging the table's
# formatting one column to the right.
$ws.Range("J3:J5").Copy($ws.Range("K3:K5"))

# New data for year 2020.
$ws.Range("K4").Value = 2020
$ws.Range("K5").Value = 173

# Restore the saved cursor/selection position.
$ws.Range("I18").Select()
